# "Actualizacion desde MV -datos-"
# Append the next 6 daily rows (09-10-2021 .. 14-10-2021) to the
# "reinversion de bonos bancarios" daily table, continuing the existing
# Serie / Bonos Bancarios en UF / Bonos Bancarios en Pesos ($) columns
# with the same values as the last existing row (449, 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated row in column A (currently row 255 -> 08-10-2021).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$dates = @("09-10-2021", "10-10-2021", "11-10-2021", "12-10-2021", "13-10-2021", "14-10-2021")

# Far-away helper cell used to push each date through Excel's formula/paste
# pipeline so it lands back in column A as plain text, instead of being
# silently auto-parsed into a date serial (which happens for ambiguous
# day-of-month values <= 12, e.g. 09, 10, 11, 12).
$helper = $ws.Cells.Item($ws.Rows.Count, $ws.Columns.Count)

$r = $lastRow
foreach ($d in $dates) {
    $r = $r + 1
    $target = $ws.Cells.Item($r, 1)

    # A trailing space defeats the date auto-detection heuristic, so this
    # literal is stored verbatim as text.
    $target.Value = $d + " "

    # Trim the trailing space back off via a formula evaluated elsewhere,
    # then paste the resulting text *value only* back into the cell -
    # PasteSpecial(values) does not re-run date auto-detection.
    $helper.Formula = "=LEFT(" + $target.Address($false, $false) + ",10)"
    $helper.Copy()
    $target.PasteSpecial(-4163)

    $ws.Cells.Item($r, 2).Value = 449
    $ws.Cells.Item($r, 3).Value = 0
}

$helper.Clear()
$excel.CutCopyMode = 0
